# AGA206 Assessment 2 Checklist - "SpeedRunMode and Free camera"
#
# 1. Rename the 22nd optional-module task (row 35) from "Come up with your
#    own" to "NPC with dialogue" (a brand-new entry in the shared strings
#    table, leaving the existing "Come up with your own" string used by
#    row 36 untouched).
# 2. Tick all four mandatory-module checkboxes (J12:J15), which flips their
#    linked "Done"/"To Be Done" status formulas and the mandatory-modules-
#    completed counter in D7.
# 3. Leave the cursor on D35 (matching where the author was last working).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename task 22 -------------------------------------------------------
$ws.Range("C35").Value = "NPC with dialogue"

# --- Check the four mandatory-module boxes --------------------------------
$ws.Range("J12").Value = $true
$ws.Range("J13").Value = $true
$ws.Range("J14").Value = $true
$ws.Range("J15").Value = $true

# --- Match the author's final selection ------------------------------------
$ws.Range("D35").Select()
